$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, pushing the existing rows 44-45 down to 45-46.
$ws.Rows("44:44").Insert()

# Populate the newly inserted row 44 with the new weekly record.
$ws.Range("A44").Value = 4
$ws.Range("B44").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C44").Value = "Los Lagos"
$ws.Range("D44").Value = 44946
$ws.Range("E44").Value = 10
$ws.Range("F44").Value = "Fruta"
$ws.Range("G44").Value = 100101
$ws.Range("H44").Value = "Berries"
$ws.Range("I44").Value = 100101001
$ws.Range("J44").Value = "Arándano (blue)"
$ws.Range("K44").Value = "Sin especificar"
$ws.Range("L44").Value = "Primera"
$ws.Range("M44").Value = 240
$ws.Range("N44").Value = 2000
$ws.Range("O44").Value = 2200
$ws.Range("P44").Value = 2100
$ws.Range("Q44").Value = "$/bandeja 2 kilos"
$ws.Range("R44").Value = "Provincia de Curicó"
$ws.Range("S44").Value = 1050
$ws.Range("T44").Value = 2
